$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + report date range) ---
$ws.Range("A8").Value = "Volume 29   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/7/2022  Through  11/13/2022"

# --- Crime statistics table updates (rows 15-29) ---
$ws.Range("M15").Value = -23.529411764705
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -81.818181818181
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -14.285714285714
$ws.Range("I16").Value = 179
$ws.Range("J16").Value = 207
$ws.Range("K16").Value = -13.526570048309
$ws.Range("L16").Value = 23.448275862069
$ws.Range("M16").Value = -14.354066985645
$ws.Range("N16").Value = 141.891891891892
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -42.857142857142
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 46.666666666666
$ws.Range("I17").Value = 242
$ws.Range("J17").Value = 271
$ws.Range("K17").Value = -10.701107011070
$ws.Range("L17").Value = 33.701657458563
$ws.Range("M17").Value = 49.382716049382
$ws.Range("N17").Value = 340
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -55.555555555555
$ws.Range("I18").Value = 84
$ws.Range("J18").Value = 115
$ws.Range("K18").Value = -26.956521739130
$ws.Range("L18").Value = -59.615384615384
$ws.Range("M18").Value = -8.695652173913
$ws.Range("N18").Value = -2.325581395348
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -37.5
$ws.Range("I19").Value = 259
$ws.Range("J19").Value = 306
$ws.Range("K19").Value = -15.359477124183
$ws.Range("L19").Value = 18.264840182648
$ws.Range("M19").Value = 10.212765957446
$ws.Range("N19").Value = 516.666666666667
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -72.727272727272
$ws.Range("I20").Value = 120
$ws.Range("J20").Value = 66
$ws.Range("K20").Value = 81.818181818181
$ws.Range("L20").Value = 62.162162162162
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 69.014084507042
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -52
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -23.595505617977
$ws.Range("I21").Value = 899
$ws.Range("J21").Value = 980
$ws.Range("K21").Value = -8.265306122448
$ws.Range("L21").Value = 6.896551724137
$ws.Range("M21").Value = 15.404364569961
$ws.Range("N21").Value = 169.161676646707
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("D23").Value = "0"
$ws.Range("E23").Value = "***.*"
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -43.75
$ws.Range("G24").Value = 162
$ws.Range("H24").Value = -54.938271604938
$ws.Range("I24").Value = 1168
$ws.Range("J24").Value = 1367
$ws.Range("K24").Value = -14.557425018288
$ws.Range("L24").Value = 88.083735909822
$ws.Range("M24").Value = 127.237354085603
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -22.222222222222
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = -6.060606060606
$ws.Range("I25").Value = 355
$ws.Range("J25").Value = 367
$ws.Range("K25").Value = -3.269754768392
$ws.Range("L25").Value = 39.215686274509
$ws.Range("M25").Value = -16.075650118203
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -25
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = -4.347826086956
$ws.Range("L27").Value = 69.230769230769
$ws.Range("G28").Value = "0"
$ws.Range("H28").Value = "***.*"
$ws.Range("N28").Value = -25
$ws.Range("G29").Value = "0"
$ws.Range("H29").Value = "***.*"
$ws.Range("N29").Value = -33.333333333333
